{"js": "// Ultimi controlli ai deliverables\n// The functional-requirements bullet about seat reservations is being\n// trimmed: it should only talk about \"tables\" (not \"tables and\n// tournaments\"), and the timing clause should only mention the 1-hour/\n// table rule (the \"3 hours from the start of the tournament\" half is\n// dropped).\n//\n// Original sentence:\n//   \"The system shall provide seat reservations at tables and\n//    tournaments only to users with player-type accounts and only\n//    within 1 hour from the start of the table and 3 hours from the\n//    start of the tournament.\"\n// New sentence:\n//   \"The system shall provide seat reservations at tables only to\n//    users with player-type accounts and only within 1 hour from the\n//    start of the table.\"\n//\n// Both deleted spans are unique substrings in the document, so a plain\n// body-wide search/replace (replace with \"\") is enough and is robust to\n// however the original text happens to be split across <w:r> runs.\n\nconst body = context.document.body;\n\nconst firstHit = body.search(\"s and tournament\", { matchCase: true, matchWholeWord: false });\nfirstHit.load(\"items\");\nawait context.sync();\n\nif (firstHit.items.length === 0) {\n  throw new Error(\"Could not find the text 's and tournament' to remove.\");\n}\nfirstHit.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n\nconst secondHit = body.search(\" and 3 hours from the start of the tournament\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsecondHit.load(\"items\");\nawait context.sync();\n\nif (secondHit.items.length === 0) {\n  throw new Error(\"Could not find the text ' and 3 hours from the start of the tournament' to remove.\");\n}\nsecondHit.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n", "ps1": "# Ultimi controlli ai deliverables\n#\n# The functional-requirements bullet about seat reservations is being\n# trimmed: it should only talk about \"tables\" (not \"tables and\n# tournaments\"), and the timing clause should only mention the 1-hour/\n# table rule (the \"3 hours from the start of the tournament\" half is\n# dropped).\n#\n# Original sentence:\n#   \"The system shall provide seat reservations at tables and\n#    tournaments only to users with player-type accounts and only\n#    within 1 hour from the start of the table and 3 hours from the\n#    start of the tournament.\"\n# New sentence:\n#   \"The system shall provide seat reservations at tables only to\n#    users with player-type accounts and only within 1 hour from the\n#    start of the table.\"\n#\n# Both deleted spans are unique substrings in the document, so a plain\n# Find/Replace (wdReplaceAll, replace with \"\") is enough and is robust\n# to however the original text happens to be split across runs.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"s and tournament\"\n$find1.Replacement.Text = \"\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceAll)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \" and 3 hours from the start of the tournament\"\n$find2.Replacement.Text = \"\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
